$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 1: title-slide subtitle "HILT 2015" -> "HILT 2016"
# ---------------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$subTitle = $s1.Shapes.Item(2)
$tr1 = $subTitle.TextFrame.TextRange
$lastChar = $tr1.Characters($tr1.Length, 1)
$lastChar.Text = "6"

# ---------------------------------------------------------------------------
# Slide 19: "Interpreter" section -> "Running Ruby"
# ---------------------------------------------------------------------------
$s19 = $p.Slides.Item(19)

# Title shape: "Interpreter" -> "Running Ruby"
$title = $s19.Shapes.Item(1)
$trTitle = $title.TextFrame.TextRange
$trTitle.Text = "Running Ruby"
$trTitle.LanguageID = "en-US"

# Body shape: trim the explanatory bullets down to a single merged bullet
$body = $s19.Shapes.Item(2)
$trBody = $body.TextFrame.TextRange

# Remove the four paragraphs between the first bullet and the
# "different ways to run code" bullet:
#   "Its code cannot be run directly"
#   "It must be run through a Ruby interpreter"
#   "Most common interpreter is Matz's Ruby Interpreter (MRI)"
#   "There are others (jruby, rubinius, etc.)"
# (characters 33-268 of the body text span these four paragraphs plus the
# trailing "There are different ways to run code through a Ruby interpretor"
# paragraph's own text, which we also drop here and re-type below).
$midSpan = $trBody.Characters(33, (268 - 33 + 1))
$midSpan.Delete()

# Replace the first paragraph's text ("Ruby is an interpreted language")
# with the merged "There are different ways to run code through a Ruby
# interpretor" sentence, split into two runs the same way the final deck
# has it.
$firstPara = $trBody.Characters(1, 31)
$firstPara.Text = "There are different ways to run code through a Ruby interpretor"

$thereRun = $trBody.Characters(1, 6)
$thereRun.Text = "There "

# ---------------------------------------------------------------------------
# Slide 8: no textual change other than the OOXML dropping a redundant
# trailing endParaRPr on the last "[...]" paragraph -- nothing user visible
# to replicate through the object model.
# ---------------------------------------------------------------------------
